$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2254.4285
$ws.Range("I40").Value = 2829.6155
$ws.Range("J40").Value = 1319.75
$ws.Range("K40").Value = 2829.6155
$ws.Range("L40").Value = 1319.75
$ws.Range("M40").Value = -2654.6155
$ws.Range("N40").Value = -1669.75

$ws.Range("H62").Value = 19234592
$ws.Range("I62").Value = 35716316
$ws.Range("J62").Value = 5913.4165
$ws.Range("K62").Value = 35716316
$ws.Range("L62").Value = 5913.4165
$ws.Range("M62").Value = -35715692
$ws.Range("N62").Value = -7161.4165

$ws.Range("H65").Value = 19234592
$ws.Range("I65").Value = 35716316
$ws.Range("J65").Value = 5913.4165
$ws.Range("K65").Value = 178581580
$ws.Range("L65").Value = 29567.0825
$ws.Range("M65").Value = -178578460
$ws.Range("N65").Value = -35807.0825

$ws.Range("H92").Value = 3310.5557
$ws.Range("I92").Value = 3406.0667
$ws.Range("K92").Value = 3406.0667
$ws.Range("M92").Value = -2158.0667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 5333.3335
$ws.Range("J46").Value = 6000
$ws.Range("L46").Value = 6000
$ws.Range("N46").Value = -6638

$ws.Range("H55").Value = 43333.332
$ws.Range("J55").Value = 65000
$ws.Range("L55").Value = 65000
$ws.Range("N55").Value = -65630

$ws.Range("H61").Value = 4266.591
$ws.Range("I61").Value = 2003.8235
$ws.Range("K61").Value = 2003.8235
$ws.Range("M61").Value = -1791.8235

$ws.Range("H102").Value = 2952.4285
$ws.Range("I102").Value = 2611.1667
$ws.Range("K102").Value = 2611.1667
$ws.Range("M102").Value = -989.1667000000002

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H122").Value = 2524.1904
$ws.Range("I122").Value = 1944.25
$ws.Range("K122").Value = 5832.75
$ws.Range("M122").Value = -3382.75

$ws.Range("H136").Value = 4266.591
$ws.Range("I136").Value = 2003.8235
$ws.Range("K136").Value = 6011.470499999999
$ws.Range("M136").Value = -3461.470499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 53000
$ws.Range("J35").Value = 53000
$ws.Range("L35").Value = 53000
$ws.Range("N35").Value = -53620

$ws.Range("H99").Value = 3888.7778
$ws.Range("J99").Value = 2992
$ws.Range("L99").Value = 2992
$ws.Range("N99").Value = -5988

$ws.Range("H132").Value = 103998.664
$ws.Range("J132").Value = 103998.664
$ws.Range("L132").Value = 103998.664
$ws.Range("N132").Value = -114118.664

$ws.Range("H134").Value = 3035.4062
$ws.Range("I134").Value = 2468.6155
$ws.Range("K134").Value = 7405.8465
$ws.Range("M134").Value = -4870.8465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 411.04
$ws.Range("I7").Value = 376.2857
$ws.Range("J7").Value = 455.27274
$ws.Range("K7").Value = 376.2857
$ws.Range("L7").Value = 455.27274
$ws.Range("M7").Value = -263.2857
$ws.Range("N7").Value = -681.27274

$ws.Range("H22").Value = 304.05713
$ws.Range("I22").Value = 309.96774
$ws.Range("J22").Value = 258.25
$ws.Range("K22").Value = 309.96774
$ws.Range("L22").Value = 258.25
$ws.Range("M22").Value = 40.03226000000001
$ws.Range("N22").Value = -958.25

$ws.Range("H31").Value = 386783.8
$ws.Range("I31").Value = 715342.3
$ws.Range("K31").Value = 715342.3
$ws.Range("M31").Value = -715047.3

$ws.Range("H34").Value = 386783.8
$ws.Range("I34").Value = 715342.3
$ws.Range("K34").Value = 715342.3
$ws.Range("M34").Value = -715140.3

$ws.Range("H60").Value = 8999
$ws.Range("I60").Value = 8999
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 8999
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -8488
$ws.Range("N60").ClearContents()

$ws.Range("H94").Value = 1287.75
$ws.Range("I94").Value = 762.1111
$ws.Range("K94").Value = 762.1111
$ws.Range("M94").Value = -311.1111

$ws.Range("H95").Value = 13794.4
$ws.Range("J95").Value = 13794.4
$ws.Range("L95").Value = 13794.4
$ws.Range("N95").Value = -19286.4

$ws.Range("H105").Value = 5607.1377
$ws.Range("I105").Value = 2062.4
$ws.Range("K105").Value = 2062.4
$ws.Range("M105").Value = -315.4000000000001

$ws.Range("H107").Value = 5618.533
$ws.Range("I107").Value = 779.4286
$ws.Range("J107").Value = 7091.304
$ws.Range("K107").Value = 779.4286
$ws.Range("L107").Value = 7091.304
$ws.Range("M107").Value = 1140.5714
$ws.Range("N107").Value = -10931.304

$ws.Range("H132").Value = 3400.08
$ws.Range("I132").Value = 3159.2727
$ws.Range("K132").Value = 9477.8181
$ws.Range("M132").Value = -6947.8181

$ws.Range("H134").Value = 5935.278
$ws.Range("I134").Value = 6231.4517
$ws.Range("K134").Value = 18694.3551
$ws.Range("M134").Value = -16159.3551

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 155.12
$ws.Range("I2").Value = 190
$ws.Range("J2").Value = 122.92308
$ws.Range("K2").Value = 1140
$ws.Range("L2").Value = 737.5384799999999
$ws.Range("M2").Value = -1027
$ws.Range("N2").Value = -963.5384799999999

$ws.Range("H5").Value = 1025.75
$ws.Range("I5").Value = 999
$ws.Range("J5").Value = 1052.5
$ws.Range("K5").Value = 2997
$ws.Range("L5").Value = 3157.5
$ws.Range("M5").Value = -2885
$ws.Range("N5").Value = -3381.5

$ws.Range("H12").Value = 276.75
$ws.Range("I12").Value = 507
$ws.Range("K12").Value = 1521
$ws.Range("M12").Value = -1348

$ws.Range("H103").Value = 2132.8572
$ws.Range("J103").Value = 2479
$ws.Range("L103").Value = 7437
$ws.Range("N103").Value = -9195

$ws.Range("H121").Value = 104707.3
$ws.Range("J121").Value = 115564
$ws.Range("L121").Value = 346692
$ws.Range("N121").Value = -349312

$ws.Range("H131").Value = 1427.4667
$ws.Range("J131").Value = 1441.2106
$ws.Range("L131").Value = 4323.6318
$ws.Range("N131").Value = -14403.6318

$ws.Range("H135").Value = 1025.75
$ws.Range("I135").Value = 999
$ws.Range("J135").Value = 1052.5
$ws.Range("K135").Value = 8991
$ws.Range("L135").Value = 9472.5
$ws.Range("M135").Value = -6456
$ws.Range("N135").Value = -14542.5

$ws.Range("H140").Value = 8073011.5
$ws.Range("I140").Value = 35743620
$ws.Range("J140").Value = 2417.2083
$ws.Range("K140").Value = 107230860
$ws.Range("L140").Value = 7251.624899999999
$ws.Range("M140").Value = -107225680
$ws.Range("N140").Value = -17611.6249

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 29997.5
$ws.Range("J15").Value = 29997.5
$ws.Range("L15").Value = 29997.5
$ws.Range("N15").Value = -30573.5

$ws.Range("H70").Value = 8448.474
$ws.Range("J70").Value = 8478.846
$ws.Range("L70").Value = 8478.846
$ws.Range("N70").Value = -9018.846

$ws.Range("H73").Value = 8448.474
$ws.Range("J73").Value = 8478.846
$ws.Range("L73").Value = 8478.846
$ws.Range("N73").Value = -10350.846

$ws.Range("H81").Value = 29997.5
$ws.Range("J81").Value = 29997.5
$ws.Range("L81").Value = 29997.5
$ws.Range("N81").Value = -31993.5

$ws.Range("H84").Value = 29997.5
$ws.Range("J84").Value = 29997.5
$ws.Range("L84").Value = 89992.5
$ws.Range("N84").Value = -99976.5

$ws.Range("H126").Value = 10662.333
$ws.Range("I126").Value = 18727.75
$ws.Range("J126").Value = 4210
$ws.Range("K126").Value = 56183.25
$ws.Range("L126").Value = 12630
$ws.Range("M126").Value = -53713.25
$ws.Range("N126").Value = -17570

$ws.Range("H132").Value = 38840.066
$ws.Range("I132").Value = 47617.477
$ws.Range("K132").Value = 142852.431
$ws.Range("M132").Value = -140322.431

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4495.727
$ws.Range("I122").Value = 4411.8
$ws.Range("K122").Value = 13235.4
$ws.Range("M122").Value = -10785.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 13013
$ws.Range("J31").Value = 13013
$ws.Range("L31").Value = 13013
$ws.Range("N31").Value = -13709

$ws.Range("H41").Value = 33699.75
$ws.Range("J41").Value = 33699.75
$ws.Range("L41").Value = 33699.75
$ws.Range("N41").Value = -34479.75

$ws.Range("H74").Value = 29855.715
$ws.Range("I74").Value = 12999.5
$ws.Range("K74").Value = 12999.5
$ws.Range("M74").Value = -12063.5

$ws.Range("H77").Value = 29855.715
$ws.Range("I77").Value = 12999.5
$ws.Range("K77").Value = 38998.5
$ws.Range("M77").Value = -34318.5

$ws.Range("H81").Value = 3956.75
$ws.Range("I81").Value = 4548.9
$ws.Range("J81").Value = 996
$ws.Range("K81").Value = 9097.799999999999
$ws.Range("L81").Value = 1992
$ws.Range("M81").Value = -8036.799999999999
$ws.Range("N81").Value = -4114

$ws.Range("H84").Value = 3956.75
$ws.Range("I84").Value = 4548.9
$ws.Range("J84").Value = 996
$ws.Range("K84").Value = 45489
$ws.Range("L84").Value = 9960
$ws.Range("M84").Value = -40185
$ws.Range("N84").Value = -20568

$ws.Range("H132").Value = 1804.963
$ws.Range("I132").Value = 1730.28
$ws.Range("K132").Value = 5190.84
$ws.Range("M132").Value = -2660.84

$ws.Range("H133").Value = 35000
$ws.Range("J133").Value = 35000
$ws.Range("L133").Value = 35000
$ws.Range("N133").Value = -45120
